# Improved biplots with Rezek
# Adds replicate TEF (H) values and a STDEV.S summary (I) for the
# Phytoplankton/phyto block, notes that the other phyto block's SD was
# improvised (0.5, no variance reported in Fry 2002) and references that
# note from G38, and updates the selection to the newly edited area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New replicate TEF values in column H (rows 36-38) and their STDEV.S in I36
$ws.Range("H36").Value = 19.5
$ws.Range("H37").Value = 23
$ws.Range("H38").Value = 26
$ws.Range("I36").Formula = "=_xlfn.STDEV.S(H36:H38)"

# Note explaining the improvised SD, stored at I43 and referenced by G38
$ws.Range("I43").Value = "sd set to 0.5 no variance reported in Fry 2002 so improvised sd"
$ws.Range("G38").Formula = "=I43"

# G43 no longer computed via STDEVA - just the plain value remains
$ws.Range("G43").Value = 0

# Move the active selection to the newly added cell
$ws.Range("I44").Select()
